# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "25.790.23"
$ws.Cells.Item(2, 5).Value = "  -0.11%  "
$ws.Cells.Item(3, 4).Value = "1.636.22"
$ws.Cells.Item(3, 5).Value = "  +0.16%  "
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "'215.65"
$ws.Cells.Item(5, 5).Value = "  +0.38%  "
$ws.Cells.Item(6, 5).Value = "  -0.61%  "
$ws.Cells.Item(7, 5).Value = "  -0.10%  "
$ws.Cells.Item(8, 5).Value = "  +0.08%  "
$ws.Cells.Item(9, 4).Value = "'0.0635"
$ws.Cells.Item(9, 5).Value = "  -1.10%  "
$ws.Cells.Item(10, 4).Value = "'19.57"
$ws.Cells.Item(10, 5).Value = "  -1.54%  "
$ws.Cells.Item(11, 5).Value = "  +1.59%  "
$ws.Cells.Item(12, 5).Value = "  +0.23%  "
$ws.Cells.Item(13, 4).Value = "1.864.05"
$ws.Cells.Item(13, 5).Value = "  +0.28%  "
$ws.Cells.Item(14, 4).Value = "1.641.72"
$ws.Cells.Item(14, 5).Value = "  +0.12%  "
$ws.Cells.Item(15, 4).Value = "'0.563"
$ws.Cells.Item(15, 5).Value = "  +1.04%  "
$ws.Cells.Item(16, 4).Value = "0.0₃0763"
$ws.Cells.Item(16, 5).Value = "  -0.53%  "
$ws.Cells.Item(17, 4).Value = "'63.26"
$ws.Cells.Item(17, 5).Value = "  +0.36%  "
$ws.Cells.Item(18, 4).Value = "25.827.14"
$ws.Cells.Item(18, 5).Value = "  +0.00%  "
$ws.Cells.Item(19, 5).Value = "  -0.11%  "
$ws.Cells.Item(20, 5).Value = "  +2.26%  "
$ws.Cells.Item(21, 4).Value = "'192.58"
$ws.Cells.Item(22, 4).Value = "'9.98"
$ws.Cells.Item(22, 5).Value = "  +0.57%  "
$ws.Cells.Item(23, 4).Value = "'6.28"
$ws.Cells.Item(23, 5).Value = "  +1.73%  "
$ws.Cells.Item(24, 4).Value = "'1.83"
$ws.Cells.Item(24, 5).Value = "  +4.53%  "
$ws.Cells.Item(25, 5).Value = "  -0.09%  "
$ws.Cells.Item(26, 4).Value = "'141.48"
$ws.Cells.Item(26, 5).Value = "  +1.45%  "
$ws.Cells.Item(27, 5).Value = "  +1.63%  "
$ws.Cells.Item(28, 4).Value = "'6.91"
$ws.Cells.Item(28, 5).Value = "  +1.17%  "
$ws.Cells.Item(29, 4).Value = "'15.49"
$ws.Cells.Item(29, 5).Value = "  -0.10%  "
$ws.Cells.Item(30, 5).Value = "  +0.06%  "
$ws.Cells.Item(31, 5).Value = "  -0.50%  "
$ws.Cells.Item(32, 4).Value = "'3.32"
$ws.Cells.Item(32, 5).Value = "  +0.42%  "
$ws.Cells.Item(33, 5).Value = "  -0.50%  "
$ws.Cells.Item(34, 5).Value = "  -0.34%  "
$ws.Cells.Item(35, 5).Value = "  -0.39%  "
$ws.Cells.Item(36, 4).Value = "'0.905"
$ws.Cells.Item(36, 5).Value = "  +0.38%  "
$ws.Cells.Item(37, 4).Value = "1.134.00"
$ws.Cells.Item(37, 5).Value = "  +1.20%  "
$ws.Cells.Item(38, 5).Value = "  -1.62%  "
$ws.Cells.Item(39, 5).Value = "  -1.09%  "
$ws.Cells.Item(40, 5).Value = "  -0.28%  "
$ws.Cells.Item(41, 5).Value = "  +0.12%  "
$ws.Cells.Item(42, 2).Value = "mCoin"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Cells.Item(42, 4).Value = "'2.54"
$ws.Cells.Item(42, 5).Value = "  +0.77%  "
$ws.Cells.Item(43, 2).Value = "FraxShare"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(43, 4).Value = "'5.58"
$ws.Cells.Item(43, 5).Value = "  +0.90%  "
$ws.Cells.Item(44, 2).Value = "Quant"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(44, 4).Value = "'100.42"
$ws.Cells.Item(44, 5).Value = "  +0.85%  "
$ws.Cells.Item(45, 2).Value = "TrustWalletToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(45, 4).Value = "'0.803"
$ws.Cells.Item(45, 5).Value = "  +0.50%  "
$ws.Cells.Item(46, 2).Value = "RocketPoolETH"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Cells.Item(46, 4).Value = "1.772.56"
$ws.Cells.Item(46, 5).Value = "  -0.16%  "
$ws.Cells.Item(47, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Cells.Item(47, 4).Value = "0.0₆0112"
$ws.Cells.Item(47, 5).Value = "  +3.99%  "
$ws.Cells.Item(48, 2).Value = "Aave"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Cells.Item(48, 4).Value = "'55.34"
$ws.Cells.Item(48, 5).Value = "  -0.16%  "
$ws.Cells.Item(49, 2).Value = "Mantle"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(49, 4).Value = "'0.417"
$ws.Cells.Item(49, 5).Value = "  -1.13%  "
$ws.Cells.Item(50, 2).Value = "Cronos"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(50, 4).Value = "'0.0502"
$ws.Cells.Item(50, 5).Value = "  -0.17%  "
$ws.Cells.Item(51, 2).Value = "RenderToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(51, 4).Value = "'1.43"
$ws.Cells.Item(51, 5).Value = "  +4.05%  "

Write-Output "Applied cryptos update"
